# Scheduled data refresh: update market-price-derived columns (H:N)
# across the Leve profit tables on each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 3759.3572
$ws.Range("I41").Value = 3163.2
$ws.Range("K41").Value = 3163.2
$ws.Range("M41").Value = -2723.2

$ws.Range("H70").Value = 1669.3334
$ws.Range("I70").Value = 1574.6666
$ws.Range("J70").Value = 1716.6666
$ws.Range("K70").Value = 4723.9998
$ws.Range("L70").Value = 5149.9998
$ws.Range("M70").Value = -4453.9998
$ws.Range("N70").Value = -5689.9998

$ws.Range("H73").Value = 1669.3334
$ws.Range("I73").Value = 1574.6666
$ws.Range("J73").Value = 1716.6666
$ws.Range("K73").Value = 4723.9998
$ws.Range("L73").Value = 5149.9998
$ws.Range("M73").Value = -3787.9998
$ws.Range("N73").Value = -7021.9998

$ws.Range("H80").Value = 1616.75
$ws.Range("I80").Value = 1375.7
$ws.Range("J80").Value = 1857.8
$ws.Range("K80").Value = 4127.1
$ws.Range("L80").Value = 5573.4
$ws.Range("M80").Value = -3129.1
$ws.Range("N80").Value = -7569.4

$ws.Range("H83").Value = 1616.75
$ws.Range("I83").Value = 1375.7
$ws.Range("J83").Value = 1857.8
$ws.Range("K83").Value = 12381.3
$ws.Range("L83").Value = 16720.2
$ws.Range("M83").Value = -7389.300000000001
$ws.Range("N83").Value = -26704.2

$ws.Range("H125").Value = 850
$ws.Range("I125").Value = 700
$ws.Range("K125").Value = 6300
$ws.Range("M125").Value = -3840

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1122.5
$ws.Range("I61").Value = 1122.5
$ws.Range("K61").Value = 1122.5
$ws.Range("M61").Value = -910.5

$ws.Range("H88").Value = 1990
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 2485
$ws.Range("K88").Value = 1000
$ws.Range("L88").Value = 2485
$ws.Range("M88").Value = -594
$ws.Range("N88").Value = -3297

$ws.Range("H91").Value = 1990
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 2485
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 2485
$ws.Range("M91").Value = 404
$ws.Range("N91").Value = -5293

$ws.Range("H110").Value = 368
$ws.Range("I110").Value = 368
$ws.Range("K110").Value = 368
$ws.Range("M110").Value = 1677

$ws.Range("H122").Value = 3845.9167
$ws.Range("I122").Value = 3530
$ws.Range("J122").Value = 4288.2
$ws.Range("K122").Value = 10590
$ws.Range("L122").Value = 12864.6
$ws.Range("M122").Value = -8140
$ws.Range("N122").Value = -17764.6

$ws.Range("H136").Value = 1122.5
$ws.Range("I136").Value = 1122.5
$ws.Range("K136").Value = 3367.5
$ws.Range("M136").Value = -817.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7916.3335
$ws.Range("I107").Value = 10499.5
$ws.Range("J107").Value = 2750
$ws.Range("K107").Value = 10499.5
$ws.Range("L107").Value = 2750
$ws.Range("M107").Value = -8579.5
$ws.Range("N107").Value = -6590

$ws.Range("H134").Value = 2378.25
$ws.Range("I134").Value = 747.5
$ws.Range("K134").Value = 2242.5
$ws.Range("M134").Value = 292.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3091.125
$ws.Range("I16").Value = 2927.75
$ws.Range("K16").Value = 2927.75
$ws.Range("M16").Value = -2640.75

$ws.Range("H113").Value = 3091.125
$ws.Range("I113").Value = 2927.75
$ws.Range("K113").Value = 2927.75
$ws.Range("M113").Value = -757.75

$ws.Range("H134").Value = 6250
$ws.Range("I134").Value = 1000
$ws.Range("K134").Value = 3000
$ws.Range("M134").Value = -465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 45464.047
$ws.Range("I2").Value = 76936
$ws.Range("J2").Value = 4.5555553
$ws.Range("K2").Value = 461616
$ws.Range("L2").Value = 27.3333318
$ws.Range("M2").Value = -461503
$ws.Range("N2").Value = -253.3333318

$ws.Range("H17").Value = 41.666668
$ws.Range("J17").Value = 23.666666
$ws.Range("L17").Value = 70.99999800000001
$ws.Range("N17").Value = -408.999998

$ws.Range("H113").Value = 413.5
$ws.Range("I113").Value = 453.57144
$ws.Range("K113").Value = 1360.71432
$ws.Range("M113").Value = 809.28568

$ws.Range("H137").Value = 2000
$ws.Range("J137").Value = 2000
$ws.Range("L137").Value = 6000
$ws.Range("N137").Value = -16200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 2178315.5
$ws.Range("J7").Value = 2178315.5
$ws.Range("L7").Value = 2178315.5
$ws.Range("N7").Value = -2178539.5

$ws.Range("H8").Value = 2178315.5
$ws.Range("J8").Value = 2178315.5
$ws.Range("L8").Value = 2178315.5
$ws.Range("N8").Value = -2178593.5

$ws.Range("H9").Value = 10996.5
$ws.Range("I9").Value = 1993
$ws.Range("K9").Value = 1993
$ws.Range("M9").Value = -1823

$ws.Range("H11").Value = 6218800
$ws.Range("I11").Value = 8318254.5
$ws.Range("K11").Value = 8318254.5
$ws.Range("M11").Value = -8318115.5

$ws.Range("H13").Value = 1020.7
$ws.Range("I13").Value = 89
$ws.Range("J13").Value = 1952.4
$ws.Range("K13").Value = 89
$ws.Range("L13").Value = 1952.4
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = -2230.4

$ws.Range("H70").Value = 1971.6
$ws.Range("I70").Value = 1119.3334
$ws.Range("J70").Value = 3250
$ws.Range("K70").Value = 1119.3334
$ws.Range("L70").Value = 3250
$ws.Range("M70").Value = -849.3334
$ws.Range("N70").Value = -3790

$ws.Range("H73").Value = 1971.6
$ws.Range("I73").Value = 1119.3334
$ws.Range("J73").Value = 3250
$ws.Range("K73").Value = 1119.3334
$ws.Range("L73").Value = 3250
$ws.Range("M73").Value = -183.3334
$ws.Range("N73").Value = -5122

$ws.Range("H122").Value = 1885.4286
$ws.Range("I122").Value = 1866.3334
$ws.Range("K122").Value = 5599.0002
$ws.Range("M122").Value = -3149.0002

$ws.Range("H126").Value = 19000
$ws.Range("J126").Value = 27500
$ws.Range("L126").Value = 82500
$ws.Range("N126").Value = -87440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4835.125
$ws.Range("I7").Value = 1937.6
$ws.Range("K7").Value = 1937.6
$ws.Range("M7").Value = -1825.6

$ws.Range("H46").Value = 13991.429
$ws.Range("I46").Value = 9666.666999999999
$ws.Range("J46").Value = 17235
$ws.Range("K46").Value = 9666.666999999999
$ws.Range("L46").Value = 17235
$ws.Range("M46").Value = -9478.666999999999
$ws.Range("N46").Value = -17611

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988

$ws.Range("H126").Value = 4835.125
$ws.Range("I126").Value = 1937.6
$ws.Range("K126").Value = 5812.799999999999
$ws.Range("M126").Value = -3342.799999999999

$ws.Range("H132").Value = 936
$ws.Range("I132").Value = 936
$ws.Range("K132").Value = 2808
$ws.Range("M132").Value = -278

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10500
$ws.Range("I81").Value = 20000
$ws.Range("K81").Value = 40000
$ws.Range("M81").Value = -38939

$ws.Range("H82").Value = 49999.5
$ws.Range("J82").Value = 49999.5
$ws.Range("L82").Value = 49999.5
$ws.Range("N82").Value = -50765.5

$ws.Range("H84").Value = 10500
$ws.Range("I84").Value = 20000
$ws.Range("K84").Value = 200000
$ws.Range("M84").Value = -194696

$ws.Range("H85").Value = 49999.5
$ws.Range("J85").Value = 49999.5
$ws.Range("L85").Value = 49999.5
$ws.Range("N85").Value = -52651.5

$ws.Range("H107").Value = 2962.75
$ws.Range("I107").Value = 1511
$ws.Range("J107").Value = 3688.625
$ws.Range("K107").Value = 4533
$ws.Range("L107").Value = 11065.875
$ws.Range("M107").Value = -2613
$ws.Range("N107").Value = -14905.875

$ws.Range("H125").Value = 47381.668
$ws.Range("J125").Value = 47381.668
$ws.Range("L125").Value = 47381.668
$ws.Range("N125").Value = -57221.668

$ws.Range("H132").Value = 2400

$ws.Range("H136").Value = 4137.125
$ws.Range("I136").Value = 3819.6
$ws.Range("J136").Value = 4666.3335
$ws.Range("K136").Value = 11458.8
$ws.Range("L136").Value = 13999.0005
$ws.Range("M136").Value = -8908.799999999999
$ws.Range("N136").Value = -19099.0005
